$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix mis-encoded (mojibake) characters in the Regional Economic
#     Communities footnote (cell A103) ---
$ws.Range("A103").Value = "Regional Economic Communities:CEN-SAD = ""Community of Sahel-Saharan States"";COMESA = ""Common Market for Eastern and Southern Africa"";EAC = ""East African Community"";ECCAS = ""Economic Community of Central African States"";ECOWAS = ""Economic Community of West African States"";IGAD = ""Intergovernmental Authority on Development"";SADC = ""Southern African Development Community"";UMA = ""Arab Maghreb Union"";PALOP = ""Países Africanos de Língua Oficial Portuguesa"";ASEAN = ""Association of Southeast Asian Nations"";MERCOSUR = ""Mercado Común del Sur"".EU27 = ""European Union (27 members)"".OECD = ""Organisation for Economic Co-operation and Development""."

# --- Header row 2 grew taller (re-wrapped) ---
$ws.Rows.Item(2).RowHeight = 63.5

# --- Updated statistics for several country / aggregate rows ---
$ws.Range("J63").Value = 34.991578947368403
$ws.Range("K63").Value = 2.3294736842105301

$ws.Range("J65").Value = 35.174999999999997
$ws.Range("K65").Value = 4.95

$ws.Range("J66").Value = 37.119285714285702
$ws.Range("K66").Value = 10.0407142857143

$ws.Range("J76").Value = 37.314285714285703
$ws.Range("K76").Value = 2.1857142857142899

$ws.Range("J83").Value = 35.052873563218398
$ws.Range("K83").Value = 1.9747126436781599

$ws.Range("J87").Value = 35.908333333333303
$ws.Range("K87").Value = 5.4249999999999998

$ws.Range("J89").Value = 38.133333333333297

# Row 97 (aggregate row)
$ws.Range("C97").Value = 842549.25199999998
$ws.Range("D97").Value = 3576593.0168705201
$ws.Range("E97").Value = 1347361.3473996599
$ws.Range("F97").Value = 3.6103000928228699
$ws.Range("G97").Value = 4392.0388032063001
$ws.Range("I97").Value = 63.459178571428602
$ws.Range("J97").Value = 41.719230769230798
$ws.Range("K97").Value = 29.669230769230801
$ws.Range("L97").Value = 0.52607142857142997
$ws.Range("M97").Value = 0.35523076923077002

# Row 98 (aggregate row)
$ws.Range("C98").Value = 692226.44200000004
$ws.Range("D98").Value = 5590629.0987125896
$ws.Range("E98").Value = 1649113.1550898999
$ws.Range("F98").Value = 4.8789750798261302
$ws.Range("G98").Value = 8621.4504622752502
$ws.Range("I98").Value = 59.288333333333298
$ws.Range("J98").Value = 37.78
$ws.Range("K98").Value = 8.51
$ws.Range("L98").Value = 0.63606666666667
$ws.Range("M98").Value = 0.50253846153846005
